# Updated symbol list (price/volume refresh) - mirrors the scheduled
# GitHub Actions scraper commit. Values are written with a leading
# apostrophe so Excel keeps them as literal text (matching the original
# inlineStr cells) instead of re-parsing them as numbers/percentages.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'332.17"
$ws.Range("E2").Value = "'1.65%"
$ws.Range("D3").Value = "'45.88"
$ws.Range("E3").Value = "'4.31%"
$ws.Range("E4").Value = "'3.09%"
$ws.Range("D5").Value = "'0.08373"
$ws.Range("E5").Value = "'4.53%"
$ws.Range("E6").Value = "'2.33%"
$ws.Range("D7").Value = "'0.9870"
$ws.Range("E7").Value = "'3.77%"
$ws.Range("D8").Value = "'2.573"
$ws.Range("E8").Value = "'0.19%"
$ws.Range("D9").Value = "'0.1162"
$ws.Range("E9").Value = "'3.09%"
$ws.Range("D10").Value = "'0.1943"
$ws.Range("E10").Value = "'4.02%"
$ws.Range("D11").Value = "'10.43"
$ws.Range("E11").Value = "'-1.51%"
$ws.Range("D12").Value = "'0.1011"
$ws.Range("E12").Value = "'2.81%"
$ws.Range("D13").Value = "'0.04666"
$ws.Range("E13").Value = "'1.69%"
$ws.Range("D14").Value = "'0.1058"
$ws.Range("E14").Value = "'-0.60%"
$ws.Range("D15").Value = "'0.001285"
$ws.Range("E15").Value = "'1.43%"
$ws.Range("D16").Value = "'0.006134"
$ws.Range("E16").Value = "'4.64%"
$ws.Range("D17").Value = "'3.369"
$ws.Range("E17").Value = "'0.30%"
$ws.Range("D18").Value = "'4.481"
$ws.Range("E18").Value = "'4.25%"
$ws.Range("E19").Value = "'-3.21%"
$ws.Range("E20").Value = "'-0.50%"
$ws.Range("D21").Value = "'0.2647"
$ws.Range("E21").Value = "'4.03%"
$ws.Range("D22").Value = "'0.04203"
$ws.Range("E22").Value = "'2.81%"
$ws.Range("E23").Value = "'4.03%"
$ws.Range("D24").Value = "'0.004660"
$ws.Range("E24").Value = "'7.99%"
$ws.Range("E25").Value = "'10.51%"
$ws.Range("D26").Value = "'0.0003743"
$ws.Range("E26").Value = "'0.03%"
$ws.Range("D38").Value = "'0.02793"
$ws.Range("E38").Value = "'9.21%"
$ws.Range("D39").Value = "'0.05798"
$ws.Range("E39").Value = "'2.99%"
$ws.Range("D40").Value = "'0.007733"
$ws.Range("E40").Value = "'2.56%"
$ws.Range("E41").Value = "'3.00%"
$ws.Range("D42").Value = "'0.007284"
$ws.Range("E42").Value = "'-4.02%"
$ws.Range("E43").Value = "'-1.95%"
$ws.Range("D44").Value = "'0.009170"
$ws.Range("E44").Value = "'3.55%"
$ws.Range("D45").Value = "'0.3503"
$ws.Range("D46").Value = "'0.00007187"
$ws.Range("E46").Value = "'1.21%"
$ws.Range("E47").Value = "'0.14%"
$ws.Range("D48").Value = "'0.0005807"
$ws.Range("D49").Value = "'0.003504"
$ws.Range("E49").Value = "'12.67%"
$ws.Range("D50").Value = "'0.003501"
$ws.Range("E50").Value = "'-0.79%"
$ws.Range("E51").Value = "'0.14%"
